$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Dlk1"
$ws.Range("C2").Value = "Notch2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1167573333333333
$ws.Range("H2").Value = 0.350272
$ws.Range("I2").Value = 0.0006433000764991399
$ws.Range("J2").Value = 0.0006433000764991399
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 27.53580066666666
$ws.Range("N2").Value = 82.607402
$ws.Range("O2").Value = 0.2054887285464767
$ws.Range("P2").Value = 0.2054887285464768
$ws.Range("Q2").Value = 3.215006657038221
$ws.Range("R2").Value = 28.93505991334399
$ws.Range("S2").Value = 0.0001321909147936595
$ws.Range("T2").Value = 0.0001321909147936595

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Dlk1"
$ws.Range("C3").Value = "Notch2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1167573333333333
$ws.Range("H3").Value = 0.350272
$ws.Range("I3").Value = 0.0006433000764991399
$ws.Range("J3").Value = 0.0006433000764991399
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 27.50472933333333
$ws.Range("N3").Value = 82.514188
$ws.Range("O3").Value = 0.2052568555438283
$ws.Range("P3").Value = 0.2052568555438283
$ws.Range("Q3").Value = 3.211378851015111
$ws.Range("R3").Value = 28.902409659136
$ws.Range("S3").Value = 0.0001320417508733176
$ws.Range("T3").Value = 0.0001320417508733177

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Dlk1"
$ws.Range("C4").Value = "Notch2"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1167573333333333
$ws.Range("H4").Value = 0.350272
$ws.Range("I4").Value = 0.0006433000764991399
$ws.Range("J4").Value = 0.0006433000764991399
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 45.12975566666668
$ws.Range("N4").Value = 135.389267
$ws.Range("O4").Value = 0.3367854170582615
$ws.Range("P4").Value = 0.3367854170582616
$ws.Range("Q4").Value = 5.26922992562489
$ws.Range("R4").Value = 47.423069330624
$ws.Range("S4").Value = 0.0002166540845573744
$ws.Range("T4").Value = 0.0002166540845573744

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Dlk1"
$ws.Range("C5").Value = "Notch2"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1167573333333333
$ws.Range("H5").Value = 0.350272
$ws.Range("I5").Value = 0.0006433000764991399
$ws.Range("J5").Value = 0.0006433000764991399
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 33.831228
$ws.Range("N5").Value = 101.493684
$ws.Range("O5").Value = 0.2524689988514334
$ws.Range("P5").Value = 0.2524689988514334
$ws.Range("Q5").Value = 3.950043964672
$ws.Range("R5").Value = 35.550395682048
$ws.Range("S5").Value = 0.0001624133262747884
$ws.Range("T5").Value = 0.0001624133262747884

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Dlk1"
$ws.Range("C6").Value = "Notch2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 174.184255
$ws.Range("H6").Value = 522.552765
$ws.Range("I6").Value = 0.9597062674131449
$ws.Range("J6").Value = 0.959706267413145
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 27.53580066666666
$ws.Range("N6").Value = 82.607402
$ws.Range("O6").Value = 0.2054887285464767
$ws.Range("P6").Value = 0.2054887285464768
$ws.Range("Q6").Value = 4796.302924951836
$ws.Range("R6").Value = 43166.72632456653
$ws.Range("S6").Value = 0.1972088206688122
$ws.Range("T6").Value = 0.1972088206688122

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Dlk1"
$ws.Range("C7").Value = "Notch2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 174.184255
$ws.Range("H7").Value = 522.552765
$ws.Range("I7").Value = 0.9597062674131449
$ws.Range("J7").Value = 0.959706267413145
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 27.50472933333333
$ws.Range("N7").Value = 82.514188
$ws.Range("O7").Value = 0.2052568555438283
$ws.Range("P7").Value = 0.2052568555438283
$ws.Range("Q7").Value = 4790.890787903313
$ws.Range("R7").Value = 43118.01709112983
$ws.Range("S7").Value = 0.1969862906949265
$ws.Range("T7").Value = 0.1969862906949266

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Dlk1"
$ws.Range("C8").Value = "Notch2"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 174.184255
$ws.Range("H8").Value = 522.552765
$ws.Range("I8").Value = 0.9597062674131449
$ws.Range("J8").Value = 0.959706267413145
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 45.12975566666668
$ws.Range("N8").Value = 135.389267
$ws.Range("O8").Value = 0.3367854170582615
$ws.Range("P8").Value = 0.3367854170582616
$ws.Range("Q8").Value = 7860.892869130364
$ws.Range("R8").Value = 70748.03582217326
$ws.Range("S8").Value = 0.3232150755241635
$ws.Range("T8").Value = 0.3232150755241636

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Dlk1"
$ws.Range("C9").Value = "Notch2"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 174.184255
$ws.Range("H9").Value = 522.552765
$ws.Range("I9").Value = 0.9597062674131449
$ws.Range("J9").Value = 0.959706267413145
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 33.831228
$ws.Range("N9").Value = 101.493684
$ws.Range("O9").Value = 0.2524689988514334
$ws.Range("P9").Value = 0.2524689988514334
$ws.Range("Q9").Value = 5892.867244915141
$ws.Range("R9").Value = 53035.80520423626
$ws.Range("S9").Value = 0.2422960805252427
$ws.Range("T9").Value = 0.2422960805252427

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Dlk1"
$ws.Range("C10").Value = "Notch2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 7.196452999999999
$ws.Range("H10").Value = 21.589359
$ws.Range("I10").Value = 0.03965043251035593
$ws.Range("J10").Value = 0.03965043251035594
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 27.53580066666666
$ws.Range("N10").Value = 82.607402
$ws.Range("O10").Value = 0.2054887285464767
$ws.Range("P10").Value = 0.2054887285464768
$ws.Range("Q10").Value = 198.1600953150353
$ws.Range("R10").Value = 1783.440857835318
$ws.Range("S10").Value = 0.008147716962870926
$ws.Range("T10").Value = 0.008147716962870928

$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Dlk1"
$ws.Range("C11").Value = "Notch2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 7.196452999999999
$ws.Range("H11").Value = 21.589359
$ws.Range("I11").Value = 0.03965043251035593
$ws.Range("J11").Value = 0.03965043251035594
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 27.50472933333333
$ws.Range("N11").Value = 82.514188
$ws.Range("O11").Value = 0.2052568555438283
$ws.Range("P11").Value = 0.2052568555438283
$ws.Range("Q11").Value = 197.9364919250546
$ws.Range("R11").Value = 1781.428427325492
$ws.Range("S11").Value = 0.00813852309802844
$ws.Range("T11").Value = 0.008138523098028443

$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Dlk1"
$ws.Range("C12").Value = "Notch2"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 7.196452999999999
$ws.Range("H12").Value = 21.589359
$ws.Range("I12").Value = 0.03965043251035593
$ws.Range("J12").Value = 0.03965043251035594
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 45.12975566666668
$ws.Range("N12").Value = 135.389267
$ws.Range("O12").Value = 0.3367854170582615
$ws.Range("P12").Value = 0.3367854170582616
$ws.Range("Q12").Value = 324.7741655566504
$ws.Range("R12").Value = 2922.967490009853
$ws.Range("S12").Value = 0.01335368744954067
$ws.Range("T12").Value = 0.01335368744954068

$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Dlk1"
$ws.Range("C13").Value = "Notch2"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 7.196452999999999
$ws.Range("H13").Value = 21.589359
$ws.Range("I13").Value = 0.03965043251035593
$ws.Range("J13").Value = 0.03965043251035594
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 33.831228
$ws.Range("N13").Value = 101.493684
$ws.Range("O13").Value = 0.2524689988514334
$ws.Range("P13").Value = 0.2524689988514334
$ws.Range("Q13").Value = 243.464842234284
$ws.Range("R13").Value = 2191.183580108556
$ws.Range("S13").Value = 0.01001050499991589
$ws.Range("T13").Value = 0.01001050499991589
